$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.968.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.747.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5057"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -8.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2752"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06192"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07263"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.743.42"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6552"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.664"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9992"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.977.45"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006854"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.967.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.455"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.736"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.400"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.510"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.784"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.85"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.881"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08204"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.651"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04688"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.654"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9978"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6179"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.748"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01616"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.926"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9994"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.96"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3930"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7641"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.007"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1151"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.342"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.89"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05302"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3447"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.553"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.17%  "
